$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell D1 ("Divisão"), styled like the existing B1/C1 headers.
$ws.Range("D1").Value = "Divisão"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats - copy formatting only

# New data column D2:D18 - "Divisão" each team/season competed in.
$divisions = @{
    2  = "Primeira Divisão"
    3  = "Primeira Divisão"
    4  = "Primeira Divisão"
    5  = "Primeira Divisão"
    6  = "Primeira Divisão"
    7  = "Primeira Divisão"
    8  = "Primeira Divisão"
    9  = "Primeira Divisão"
    10 = "Primeira Divisão"
    11 = "Segunda Divisão"
    12 = "Primeira Divisão"
    13 = "Primeira Divisão"
    14 = "Primeira Divisão"
    15 = "Primeira Divisão"
    16 = "Primeira Divisão"
    17 = "Primeira Divisão"
    18 = "Primeira Divisão"
}

foreach ($row in $divisions.Keys) {
    $ws.Cells.Item($row, 4).Value = $divisions[$row]
}
